# Insert a new weekly price-report row for "Feria Lagunitas de Puerto Montt - Frutilla"
# at row 63, pushing the existing rows 63-149 down to 64-150 (dimension grows from
# A1:T149 to A1:T150).
#
# The new row keeps the same Mercado/Region/Producto/Categoria/Variedad/Calidad/Kg-unidad
# metadata as the row that used to occupy row 63 (it is duplicated, just like Excel does
# when you insert a row and copy the values down), while the date and price columns get
# the new week's figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row above the current row 63; everything below shifts down one row.
$ws.Rows.Item(63).Insert()

# Seed the new row 63 with the same row data that is now sitting in row 64 (i.e. what
# used to be row 63), so formatting/metadata columns start out identical.
$ws.Range("A64:T64").Copy($ws.Range("A63:T63"))

# Now overwrite the columns that actually carry the new week's values.
$ws.Range("D63").Value = 44540
$ws.Range("M63").Value = 1200
$ws.Range("N63").Value = 10000
$ws.Range("O63").Value = 10000
$ws.Range("P63").Value = 10000
$ws.Range("Q63").Value = "$/bandeja 7 kilos"
$ws.Range("R63").Value = "Provincia de Melipilla"
$ws.Range("S63").Value = 1429
